# Daily auto push: 2026-01-10 13:37 UTC
#
# The data table on Sheet1 (A:D = date / weekday / time / ranking) grows by
# two new observations for "2026/01/10" (continuing the existing 6, 10, 12,
# 15 sequence already present for that date). The two new rows are inserted
# right before the first "2026/12/29" row (row 621), which pushes every
# following row down by two (621-662 -> 623-664) and extends the used range
# from D662 to D664.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 621 by copying the formatting/type of the row
# immediately above (row 620, itself a "2026/01/10" row) and using
# Insert() on the copied row. This is the natural "insert a row like the
# one above it" gesture and keeps the new cells' styles/types (text dates,
# plain numbers) identical to their neighbours instead of letting Excel
# reinterpret a freshly-typed "2026/01/10" string as a date serial value.
$ws.Rows.Item(620).Copy()
$ws.Rows.Item(621).Insert()
$ws.Rows.Item(620).Copy()
$ws.Rows.Item(621).Insert()

# Row 621: 2026/01/10, 土, 18:00, ranking 201 (date/weekday already correct
# from the copy above, so only the time/ranking columns need updating)
$ws.Range("C621").Value = 18
$ws.Range("D621").Value = 201

# Row 622: 2026/01/10, 土, 20:00, ranking 201
$ws.Range("C622").Value = 20
$ws.Range("D622").Value = 201
